$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unnecessary "nome" column (column B), shifting the
# remaining columns left and updating references/dimension accordingly.
$ws.Columns.Item(2).Delete()

# Reflect the new active cell selection after the column removal.
$ws.Range("B1").Select()
